$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H8").Value = "CY12 1234 4321 1234 4321 1234 4321"
$ws.Range("H9").Value = "CY11 3333 2222 2222 2222 4444 3332"
$ws.Range("H10").Value = "CY33 1111 2222 2222 2222 2222 2222"
$ws.Range("H11").Value = "CY45 3333 3333 0000 3333 3333 3333"
$ws.Range("H12").Value = "IT44 5555 5555 5555 5555 5555 5555"
